$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2489.5334
$ws.Range("I40").Value = 968.4
$ws.Range("J40").Value = 3250.1
$ws.Range("K40").Value = 968.4
$ws.Range("L40").Value = 3250.1
$ws.Range("M40").Value = -793.4
$ws.Range("N40").Value = -3600.1

$ws.Range("H62").Value = 10266.667
$ws.Range("J62").Value = 11722.444
$ws.Range("L62").Value = 11722.444
$ws.Range("N62").Value = -12970.444

$ws.Range("H65").Value = 10266.667
$ws.Range("J65").Value = 11722.444
$ws.Range("L65").Value = 58612.22
$ws.Range("N65").Value = -64852.22

$ws.Range("H86").Value = 4450.75
$ws.Range("J86").Value = 4450.75
$ws.Range("L86").Value = 4450.75
$ws.Range("N86").Value = -6696.75

$ws.Range("H89").Value = 4450.75
$ws.Range("J89").Value = 4450.75
$ws.Range("L89").Value = 22253.75
$ws.Range("N89").Value = -33485.75

$ws.Range("H100").Value = 5250.3335
$ws.Range("I100").Value = 3542.1667
$ws.Range("J100").Value = 8666.666999999999
$ws.Range("K100").Value = 3542.1667
$ws.Range("L100").Value = 8666.666999999999
$ws.Range("M100").Value = -3001.1667
$ws.Range("N100").Value = -9748.666999999999

$ws.Range("H112").Value = 1760.95
$ws.Range("I112").Value = 1494.6666
$ws.Range("J112").Value = 1875.0714
$ws.Range("K112").Value = 4483.9998
$ws.Range("L112").Value = 5625.2142
$ws.Range("M112").Value = -3375.9998
$ws.Range("N112").Value = -7841.2142

$ws.Range("H132").Value = 3800.3076
$ws.Range("J132").Value = 4553.3335
$ws.Range("L132").Value = 13660.0005
$ws.Range("N132").Value = -18720.0005

$ws.Range("H137").Value = 5308.4414
$ws.Range("I137").Value = 8624.625
$ws.Range("J137").Value = 2360.7222
$ws.Range("K137").Value = 25873.875
$ws.Range("L137").Value = 7082.1666
$ws.Range("M137").Value = -23323.875
$ws.Range("N137").Value = -12182.1666

$ws.Range("H138").Value = 3503.6191
$ws.Range("I138").Value = 2622
$ws.Range("K138").Value = 7866
$ws.Range("M138").Value = -2726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5584.579
$ws.Range("I74").Value = 5114.9287
$ws.Range("J74").Value = 6899.6
$ws.Range("K74").Value = 5114.9287
$ws.Range("L74").Value = 6899.6
$ws.Range("M74").Value = -4240.9287
$ws.Range("N74").Value = -8647.6

$ws.Range("H77").Value = 5584.579
$ws.Range("I77").Value = 5114.9287
$ws.Range("J77").Value = 6899.6
$ws.Range("K77").Value = 25574.6435
$ws.Range("L77").Value = 34498
$ws.Range("M77").Value = -21206.6435
$ws.Range("N77").Value = -43234

$ws.Range("H122").Value = 4833048
$ws.Range("I122").Value = 5850184.5
$ws.Range("K122").Value = 17550553.5
$ws.Range("M122").Value = -17548103.5

$ws.Range("H139").Value = 60853.75
$ws.Range("J139").Value = 60853.75
$ws.Range("L139").Value = 60853.75
$ws.Range("N139").Value = -71133.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2199.75
$ws.Range("J20").Value = 2199.75
$ws.Range("L20").Value = 2199.75
$ws.Range("N20").Value = -2693.75

$ws.Range("H105").Value = 2240.0588
$ws.Range("I105").Value = 2077.6667
$ws.Range("J105").Value = 2422.75
$ws.Range("K105").Value = 2077.6667
$ws.Range("L105").Value = 2422.75
$ws.Range("M105").Value = -330.6667000000002
$ws.Range("N105").Value = -5916.75

$ws.Range("H107").Value = 4120.6055
$ws.Range("I107").Value = 764.5217
$ws.Range("J107").Value = 9266.6
$ws.Range("K107").Value = 764.5217
$ws.Range("L107").Value = 9266.6
$ws.Range("M107").Value = 1155.4783
$ws.Range("N107").Value = -13106.6

$ws.Range("H134").Value = 4185.1665
$ws.Range("I134").Value = 4222.2
$ws.Range("K134").Value = 12666.6
$ws.Range("M134").Value = -10131.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3888.2778
$ws.Range("I31").Value = 1954.5454
$ws.Range("J31").Value = 6927
$ws.Range("K31").Value = 1954.5454
$ws.Range("L31").Value = 6927
$ws.Range("M31").Value = -1659.5454
$ws.Range("N31").Value = -7517

$ws.Range("H34").Value = 3888.2778
$ws.Range("I34").Value = 1954.5454
$ws.Range("J34").Value = 6927
$ws.Range("K34").Value = 1954.5454
$ws.Range("L34").Value = 6927
$ws.Range("M34").Value = -1752.5454
$ws.Range("N34").Value = -7331

$ws.Range("H99").Value = 2283
$ws.Range("I99").Value = 2425
$ws.Range("K99").Value = 2425
$ws.Range("M99").Value = -927

$ws.Range("H105").Value = 2336
$ws.Range("I105").Value = 2336
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2336
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -589
$ws.Range("N105").ClearContents()

$ws.Range("H126").Value = 2283
$ws.Range("I126").Value = 2425
$ws.Range("K126").Value = 7275
$ws.Range("M126").Value = -4805

$ws.Range("H132").Value = 2416.9333
$ws.Range("I132").Value = 2411.8462
$ws.Range("J132").Value = 2450
$ws.Range("K132").Value = 7235.5386
$ws.Range("L132").Value = 7350
$ws.Range("M132").Value = -4705.5386
$ws.Range("N132").Value = -12410

$ws.Range("H134").Value = 2243
$ws.Range("I134").Value = 2216
$ws.Range("J134").Value = 2285.4285
$ws.Range("K134").Value = 6648
$ws.Range("L134").Value = 6856.2855
$ws.Range("M134").Value = -4113
$ws.Range("N134").Value = -11926.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 29530274
$ws.Range("J4").Value = 140125020
$ws.Range("L4").Value = 420375060
$ws.Range("N4").Value = -420375284

$ws.Range("H51").Value = 1128.1428
$ws.Range("I51").Value = 499
$ws.Range("K51").Value = 1497
$ws.Range("M51").Value = -1037

$ws.Range("H88").Value = 8999.799999999999
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15856

$ws.Range("H91").Value = 8999.799999999999
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17964

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9018089
$ws.Range("J70").Value = 9069.906000000001
$ws.Range("L70").Value = 9069.906000000001
$ws.Range("N70").Value = -9609.906000000001

$ws.Range("H73").Value = 9018089
$ws.Range("J73").Value = 9069.906000000001
$ws.Range("L73").Value = 9069.906000000001
$ws.Range("N73").Value = -10941.906

$ws.Range("H102").Value = 950.8
$ws.Range("I102").Value = 713.5833
$ws.Range("K102").Value = 713.5833
$ws.Range("M102").Value = 908.4167

$ws.Range("H113").Value = 47627840
$ws.Range("I113").Value = 250001920
$ws.Range("J113").Value = 10411.059
$ws.Range("K113").Value = 250001920
$ws.Range("L113").Value = 10411.059
$ws.Range("M113").Value = -249999750
$ws.Range("N113").Value = -14751.059

$ws.Range("H126").Value = 2760
$ws.Range("I126").Value = 3075
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 9225
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -6755
$ws.Range("N126").Value = -9440

$ws.Range("H132").Value = 3774.4
$ws.Range("I132").Value = 3470.5386
$ws.Range("K132").Value = 10411.6158
$ws.Range("M132").Value = -7881.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4142.7666
$ws.Range("I22").Value = 3629.7
$ws.Range("J22").Value = 4399.3
$ws.Range("K22").Value = 3629.7
$ws.Range("L22").Value = 4399.3
$ws.Range("M22").Value = -3334.7
$ws.Range("N22").Value = -4989.3

$ws.Range("H27").Value = 4142.7666
$ws.Range("I27").Value = 3629.7
$ws.Range("J27").Value = 4399.3
$ws.Range("K27").Value = 3629.7
$ws.Range("L27").Value = 4399.3
$ws.Range("M27").Value = -3522.7
$ws.Range("N27").Value = -4613.3

$ws.Range("H45").Value = 1213.6666
$ws.Range("I45").Value = 820.5
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 820.5
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -413.5
$ws.Range("N45").Value = -2814

$ws.Range("H46").Value = 3293.2354
$ws.Range("J46").Value = 3922
$ws.Range("L46").Value = 3922
$ws.Range("N46").Value = -4298

$ws.Range("H82").Value = 456497.66
$ws.Range("I82").Value = 972623.9
$ws.Range("J82").Value = 4887.25
$ws.Range("K82").Value = 972623.9
$ws.Range("L82").Value = 4887.25
$ws.Range("M82").Value = -972262.9
$ws.Range("N82").Value = -5609.25

$ws.Range("H85").Value = 456497.66
$ws.Range("I85").Value = 972623.9
$ws.Range("J85").Value = 4887.25
$ws.Range("K85").Value = 972623.9
$ws.Range("L85").Value = 4887.25
$ws.Range("M85").Value = -971375.9
$ws.Range("N85").Value = -7383.25

$ws.Range("H93").Value = 5424.9165
$ws.Range("I93").Value = 2962.5
$ws.Range("J93").Value = 6656.125
$ws.Range("K93").Value = 2962.5
$ws.Range("L93").Value = 6656.125
$ws.Range("M93").Value = -1714.5
$ws.Range("N93").Value = -9152.125

$ws.Range("H100").Value = 4324
$ws.Range("I100").Value = 2591.6365
$ws.Range("J100").Value = 7500
$ws.Range("K100").Value = 2591.6365
$ws.Range("L100").Value = 7500
$ws.Range("M100").Value = -2050.6365
$ws.Range("N100").Value = -8582

$ws.Range("H122").Value = 4662.45
$ws.Range("I122").Value = 4802.778
$ws.Range("J122").Value = 3399.5
$ws.Range("K122").Value = 14408.334
$ws.Range("L122").Value = 10198.5
$ws.Range("M122").Value = -11958.334
$ws.Range("N122").Value = -15098.5

$ws.Range("H132").Value = 5682.8335
$ws.Range("I132").Value = 5699.467
$ws.Range("K132").Value = 17098.401
$ws.Range("M132").Value = -14568.401

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3615.4
$ws.Range("I100").Value = 3792.3333
$ws.Range("J100").Value = 3350
$ws.Range("K100").Value = 7584.6666
$ws.Range("L100").Value = 6700
$ws.Range("M100").Value = -7043.6666
$ws.Range("N100").Value = -7782

$ws.Range("H122").Value = 5878.2104
$ws.Range("I122").Value = 3473.8333
$ws.Range("K122").Value = 10421.4999
$ws.Range("M122").Value = -7971.499899999999
